# Updated data spreadsheets; added signature
# Rebuild the "teaching" talks table with the latest speaking-engagement data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old trailing rows (12:18) - the new table only runs to row 11.
$ws.Rows("12:18").Delete()

# Header row stays the same set of column names, just rewritten so the
# shared-string table is rebuilt fresh for this revision of the sheet.
$ws.Range("A1").Value = "when"
$ws.Range("B1").Value = "where"
$ws.Range("C1").Value = "with"
$ws.Range("D1").Value = "what"
$ws.Range("E1").Value = "why"
$ws.Range("F1").Value = "type"

# Row 2
$ws.Range("A2").Value = 2021
$ws.Range("B2").Value = "Online"
$ws.Range("C2").Value = "AMSCUE"
$ws.Range("D2").Value = "Speaker"
$ws.Range("E2").Value = "`"Online Molecular Biology Labs`""
$ws.Range("F2").Value = "Education"

# Row 3
$ws.Range("A3").Value = 2021
$ws.Range("B3").Value = "Online"
$ws.Range("C3").Value = "UEA CEEC Rebellion"
$ws.Range("D3").Value = "Plenary"
$ws.Range("E3").Value = "`"Genetic Pest Management: knocking out pest species with applied genetics`""
$ws.Range("F3").Value = "Research"

# Row 4
$ws.Range("A4").Value = 2021
$ws.Range("B4").Value = "Online"
$ws.Range("C4").Value = "OCR Science Forum"
$ws.Range("D4").Value = "Speaker"
$ws.Range("E4").Value = "The impact of COVID-19, present and future"
$ws.Range("F4").Value = "Education"

# Row 5
$ws.Range("A5").Value = 2020
$ws.Range("B5").Value = "Online"
$ws.Range("C5").Value = "HUBS Bio-Summit"
$ws.Range("D5").Value = "Speaker"
$ws.Range("E5").Value = "`"Using Electronic Lab Notebooks to improve reflective practises in learning`""
$ws.Range("F5").Value = "Education"

# Row 6
$ws.Range("A6").Value = 2018
$ws.Range("B6").Value = "Vancouver"
$ws.Range("C6").Value = "Entomological Society of America"
$ws.Range("D6").Value = "Invited Speaker"
$ws.Range("E6").Value = "`"Localised gene drives for insect population control`""
$ws.Range("F6").Value = "Research"

# Row 7
$ws.Range("A7").Value = 2015
$ws.Range("B7").Value = "Vienna"
$ws.Range("C7").Value = "Society of Molecular Biology & Evolution"
$ws.Range("D7").Value = "Speaker"
$ws.Range("E7").Value = "`"The microbiome of the mediterranean fruit fly`""
$ws.Range("F7").Value = "Research"

# Row 8
$ws.Range("A8").Value = 2012
$ws.Range("B8").Value = "Ottawa"
$ws.Range("C8").Value = "Evolution"
$ws.Range("D8").Value = "Speaker"
$ws.Range("E8").Value = "`"What makes a successful male? Strategies for improved insect pest management`""
$ws.Range("F8").Value = "Research"

# Row 9
$ws.Range("A9").Value = 2020
$ws.Range("B9").Value = "Online"
$ws.Range("C9").Value = "Dry Labs Real Science"
$ws.Range("D9").Value = "Speaker"
$ws.Range("E9").Value = "`"Molecular Biology tools for Online teaching`""
$ws.Range("F9").Value = "Education"

# Row 10
$ws.Range("A10").Value = 2020
$ws.Range("B10").Value = "Online"
$ws.Range("C10").Value = "Higher Education Academy Talks"
$ws.Range("D10").Value = "Invited Speaker"
$ws.Range("E10").Value = "`"Synchronous on-line teaching in the biomedical sciences - Discovering how coronavirus PCR testing works`""
$ws.Range("F10").Value = "Education"

# Row 11
$ws.Range("A11").Value = 2014
$ws.Range("B11").Value = "Cambridge"
$ws.Range("C11").Value = "Department of Genetics"
$ws.Range("D11").Value = "Invited Speaker"
$ws.Range("E11").Value = "`"An introduction to genetic pest management`""
$ws.Range("F11").Value = "Research"

# Widen the "with" column to fit the longer venue/organisation names.
$ws.Columns.Item(3).ColumnWidth = 20.7

# Match the author's last on-sheet selection.
$ws.Range("F12").Select() | Out-Null
